$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G7").Value = 2.62
$ws.Range("I7").Value = 2.6
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 3.1
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 3.45
$ws.Range("Q7").Value = 1.65
$ws.Range("R7").Value = 1.98
$ws.Range("U7").Value = 1.5
$ws.Range("V7").Value = 2.25
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 15
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 32
$ws.Range("AA7").Value = 21
$ws.Range("AC7").Value = 11.25
$ws.Range("AE7").Value = 10.75
$ws.Range("AH7").Value = 10.75
$ws.Range("AL7").Value = 19.5
$ws.Range("AM7").Value = 22
$ws.Range("AN7").Value = 4.7
$ws.Range("AO7").Value = 14.5
$ws.Range("AP7").Value = 19.5
$ws.Range("AR7").Value = 90
$ws.Range("AT7").Value = 2.8
$ws.Range("AU7").Value = 6.2
$ws.Range("AW7").Value = 4.7
$ws.Range("AX7").Value = 13.5
$ws.Range("AZ7").Value = 55
$ws.Range("G9").Value = 1.6
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 5.5
$ws.Range("J9").Value = 2.1
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 13
$ws.Range("U9").Value = 1.67
$ws.Range("V9").Value = 2.1
$ws.Range("X9").Value = 8.5
$ws.Range("Z9").Value = 12
$ws.Range("AB9").Value = 21
$ws.Range("AD9").Value = 7.5
$ws.Range("AG9").Value = 151
$ws.Range("AH9").Value = 17
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 17
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 8
$ws.Range("AQ9").Value = 23
$ws.Range("AX9").Value = 26
$ws.Range("AZ9").Value = 81
$ws.Range("BB9").Value = 151
$ws.Range("V11").Value = 1.67
$ws.Range("V12").Value = 1.67
$ws.Range("V13").Value = 1.67